$d = $word.ActiveDocument
$src = $d.Range(104,105)
$ft = $src.FormattedText
$dst = $d.Range(3231, 3231)
$dst.FormattedText = $ft

$delTarget = $d.Range(3216,3231)
$delTarget.Text = ""
Write-Output "ctx: [$($d.Range(3190,3220).Text)]"
